$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tätigkeitsdokumentation")
$r = $ws.Cells.Item(3, 1)
Write-Output $r.Value
Write-Output $r.Text
Write-Output $r.Value2
